$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells whose new values look like plain decimal numbers need to be
# force-formatted as text first, otherwise Excel auto-converts them to numeric values
# and loses the original "display as typed" price formatting used throughout the sheet.

# Row 2
$ws.Range("D2").Value = "62.691.46"

# Row 3
$ws.Range("D3").Value = "3.270.74"
$ws.Range("E3").Value = "  -6.52%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.03"
$ws.Range("E5").Value = "  -3.36%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.62"
$ws.Range("E6").Value = "  -6.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.606"
$ws.Range("E7").Value = "  -4.96%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").Value = "3.258.39"
$ws.Range("E9").Value = "  -6.70%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.602"
$ws.Range("E10").Value = "  -4.67%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  -2.81%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "51.71"
$ws.Range("E12").Value = "  -3.89%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").Value = "  -4.87%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.72"
$ws.Range("E14").Value = "  -5.81%  "

# Row 15
$ws.Range("D15").Value = "3.810.16"
$ws.Range("E15").Value = "  -6.13%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "17.82"
$ws.Range("E16").Value = "  -3.48%  "

# Row 17
$ws.Range("E17").Value = "  -4.55%  "

# Row 18
$ws.Range("D18").Value = "3.271.04"
$ws.Range("E18").Value = "  -6.92%  "

# Row 19
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "62.738.90"
$ws.Range("E19").Value = "  -4.40%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.46"
$ws.Range("E20").Value = "  -5.62%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.956"
$ws.Range("E21").Value = "  -4.07%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "408.88"
$ws.Range("E22").Value = "  -1.82%  "

# Row 23
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.97"
$ws.Range("E23").Value = "  -2.20%  "

# Row 24
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.32"
$ws.Range("E24").Value = "  +4.79%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.22"
$ws.Range("E25").Value = "  +2.72%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "82.07"
$ws.Range("E26").Value = "  -4.70%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.42"
$ws.Range("E27").Value = "  -3.48%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.68"
$ws.Range("E28").Value = "  -6.20%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.41"
$ws.Range("E29").Value = "  -7.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "28.59"
$ws.Range("E30").Value = "  -5.82%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.22"
$ws.Range("E31").Value = "  -4.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "569.78"
$ws.Range("E32").Value = "  -6.52%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.17"
$ws.Range("E33").Value = "  -4.51%  "

# Row 34
$ws.Range("E34").Value = "  -4.83%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.40"
$ws.Range("E35").Value = "  -3.94%  "

# Row 36
$ws.Range("E36").Value = "  +0.50%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.143"
$ws.Range("E37").Value = "  -2.47%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "34.57"

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.35"
$ws.Range("E39").Value = "  +3.15%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0722"
$ws.Range("E40").Value = "  -8.94%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.03%  "

# Row 42
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.358"
$ws.Range("E42").Value = "  -5.82%  "

# Row 43
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.072.27"
$ws.Range("E43").Value = "  -9.90%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.22"
$ws.Range("E44").Value = "  -0.93%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.70"
$ws.Range("E45").Value = "  -5.44%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0394"
$ws.Range("E46").Value = "  -5.16%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.37"
$ws.Range("E47").Value = "  -6.57%  "

# Row 48
$ws.Range("E48").Value = "  -4.67%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.126"
$ws.Range("E49").Value = "  -4.53%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.73"
$ws.Range("E50").Value = "  -4.28%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.90"
$ws.Range("E51").Value = "  -6.49%  "
